$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-All "1234568" "2409865"
Replace-All "09.02.2024" "05.04.2024"
Replace-All "Первичная аттестация сваршиков - 2 чел." "Внеочередное освидетельствование в связи с расширением сферы деятельности ССП № 24.42.03.00765.121 от 02.05.2021"
Replace-All "Свидетельство ф. 7.1.30 № 24.42.03.00414.121 от 29.04.2024" "Свидетельство ф. 7.1.27 № 24.02.42.00987.121 от 05.05.2024"
Replace-All "536 112,20 p. (пятьсот тридцать шесть тысяч сто двенадцать рублей 20 копеек)" "15 600,00 p. (пятнадцать тысяч шестьсот рублей 00 копеек)"
Replace-All "107 222,44 p. (сто семь тысяч двести двадцать два рубля 44 копейки)" "3 120,00 p. (три тысячи сто двадцать рублей 00 копеек)"
Replace-All "643 334,64 p. (шестьсот сорок три тысячи триста тридцать четыре рубля 64 копейки)" "18 720,00 p. (восемнадцать тысяч семьсот двадцать рублей 00 копеек)"
